$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '73.373.54'
$ws.Range("E2").Value = '  +5.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.662.92'
$ws.Range("E3").Value = '  +5.74%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.11'
$ws.Range("E5").Value = '  +2.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.34'
$ws.Range("E6").Value = '  +1.93%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.530'
$ws.Range("E8").Value = '  +2.25%  '

$ws.Range("E9").Value = '  +11.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.662.38'
$ws.Range("E10").Value = '  +5.78%  '

$ws.Range("E11").Value = '  +1.04%  '

$ws.Range("E12").Value = '  +4.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.10'
$ws.Range("E13").Value = '  +1.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000192'
$ws.Range("E14").Value = '  +8.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.136.06'
$ws.Range("E15").Value = '  +6.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '73.253.32'
$ws.Range("E16").Value = '  +5.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.00'
$ws.Range("E17").Value = '  +4.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.661.25'
$ws.Range("E18").Value = '  +6.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '387.67'
$ws.Range("E19").Value = '  +6.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.64'
$ws.Range("E20").Value = '  +5.45%  '

$ws.Range("E21").Value = '  +4.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.23'
$ws.Range("E22").Value = '  +4.43%  '

$ws.Range("E23").Value = '  +22.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.60'
$ws.Range("E24").Value = '  +4.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.44'
$ws.Range("E25").Value = '  +5.42%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.95'
$ws.Range("E27").Value = '  +10.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.803.21'
$ws.Range("E28").Value = '  +6.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0976'
$ws.Range("E30").Value = '  +9.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '537.75'
$ws.Range("E31").Value = '  +5.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.11'
$ws.Range("E32").Value = '  +4.62%  '

$ws.Range("E33").Value = '  +9.60%  '

$ws.Range("E34").Value = '  +4.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.15'
$ws.Range("E36").Value = '  +1.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.43'
$ws.Range("E37").Value = '  +3.57%  '

$ws.Range("E38").Value = '  -4.00%  '

$ws.Range("E39").Value = '  +8.61%  '

$ws.Range("E40").Value = '  +2.28%  '

$ws.Range("E41").Value = '  +7.91%  '

$ws.Range("E42").Value = '  +7.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.66'
$ws.Range("E43").Value = '  +14.85%  '

$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.337'
$ws.Range("E45").Value = '  +5.10%  '

$ws.Range("E46").Value = '  +2.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '152.69'
$ws.Range("E47").Value = '  +1.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.72'
$ws.Range("E48").Value = '  +4.12%  '

$ws.Range("E49").Value = '  +6.38%  '

$ws.Range("E50").Value = '  +9.46%  '

$ws.Range("E51").Value = '  +9.12%  '
